$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 20, shifting existing rows 20..188 down to 21..189
$ws.Rows(20).Insert()

# Populate the newly inserted row 20 with the new weekly record
$ws.Cells.Item(20, 1).Value = 8
$ws.Cells.Item(20, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(20, 3).Value = "Coquimbo"
$ws.Cells.Item(20, 4).Value = 44532
$ws.Cells.Item(20, 5).Value = 4
$ws.Cells.Item(20, 6).Value = 100112012
$ws.Cells.Item(20, 7).Value = "Espinaca"
$ws.Cells.Item(20, 8).Value = "Sin especificar"
$ws.Cells.Item(20, 9).Value = "Primera"
$ws.Cells.Item(20, 10).Value = 3000
$ws.Cells.Item(20, 11).Value = 400
$ws.Cells.Item(20, 12).Value = 500
$ws.Cells.Item(20, 13).Value = 450
$ws.Cells.Item(20, 14).Value = "`$/atado 300 a 500 gramos"
$ws.Cells.Item(20, 15).Value = "Provincia del Elqu" + [char]0xED
$ws.Cells.Item(20, 16).Value = 900
$ws.Cells.Item(20, 17).Value = 0.5
$ws.Cells.Item(20, 18).Value = "Hortaliza"
